$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 updates for Row 3 (R)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 435
$wsOff.Range("C3").Value = 325
$wsOff.Range("D3").Value = 106
$wsOff.Range("E3").Value = 53

# DEF sheet - Week 17 updates for Row 3 (R)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 437
$wsDef.Range("C3").Value = 304
$wsDef.Range("D3").Value = 101
$wsDef.Range("E3").Value = 49
$wsDef.Range("F3").Value = 9
